# Update test case 20160331 - 001
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column A: timestamp strings ---
$ws.Range("A2").Value  = "20160403_220347"
$ws.Range("A3").Value  = "20160403_220619"
$ws.Range("A4").Value  = "20160403_225031"
$ws.Range("A5").Value  = "20160403_225410"
$ws.Range("A6").Value  = "20160403_225709"
$ws.Range("A7").Value  = "20160404_024810"
$ws.Range("A8").Value  = "20160404_025448"
$ws.Range("A9").Value  = "20160404_030648"
$ws.Range("A10").Value = "20160404_032232"
$ws.Range("A11").Value = "20160404_033731"

# --- Column C: Preprocess text (shared across rows 2-6 and 7-11) ---
$preprocess1 = 'convert unicode to ascii, convert to lower, trim "space" and ",", remove multiple spaces'
$preprocess2 = 'trim "space" and ",", convert to lower, remove multiple spaces, convert unicode to ascii'

$ws.Range("C2").Value = $preprocess1
$ws.Range("C3").Value = $preprocess1
$ws.Range("C4").Value = $preprocess1
$ws.Range("C5").Value = $preprocess1
$ws.Range("C6").Value = $preprocess1

$ws.Range("C7").Value = $preprocess2
$ws.Range("C8").Value = $preprocess2
$ws.Range("C9").Value = $preprocess2
$ws.Range("C10").Value = $preprocess2
$ws.Range("C11").Value = $preprocess2

# --- Column D: Features text (shared across rows 2-11) ---
$features = '4 features: #ascii/(#ascii+#digit+#punctuation), %kwAddress, %kwPhone, #max_digit_skip_0_1'
$ws.Range("D2").Value = $features
$ws.Range("D3").Value = $features
$ws.Range("D4").Value = $features
$ws.Range("D5").Value = $features
$ws.Range("D6").Value = $features
$ws.Range("D7").Value = $features
$ws.Range("D8").Value = $features
$ws.Range("D9").Value = $features
$ws.Range("D10").Value = $features
$ws.Range("D11").Value = $features

# --- Column B: RunningTime(s) ---
$ws.Range("B2").Value  = 152.133
$ws.Range("B3").Value  = 2651.614
$ws.Range("B4").Value  = 219.519
$ws.Range("B5").Value  = 179.266
$ws.Range("B6").Value  = 171.497
$ws.Range("B7").Value  = 398.191
$ws.Range("B8").Value  = 719.814
$ws.Range("B9").Value  = 943.872
$ws.Range("B10").Value = 899.649
$ws.Range("B11").Value = 560.731

# --- Column G: Test_Accuracy ---
$ws.Range("G2").Value  = 0.967333333333333
$ws.Range("G4").Value  = 0.966666666666667
$ws.Range("G5").Value  = 0.958
$ws.Range("G6").Value  = 0.957333333333333
$ws.Range("G7").Value  = 0.962
$ws.Range("G8").Value  = 0.966666666666667
$ws.Range("G9").Value  = 0.974
$ws.Range("G10").Value = 0.971333333333333
$ws.Range("G11").Value = 0.972666666666667

# --- Column J: Val_Accuracy ---
$ws.Range("J2").Value  = 0.122448979591837
$ws.Range("J3").Value  = 0.122448979591837
$ws.Range("J4").Value  = 0.122448979591837
$ws.Range("J5").Value  = 0.122448979591837
$ws.Range("J6").Value  = 0.122448979591837
$ws.Range("J7").Value  = 0.13265306122449
$ws.Range("J8").Value  = 0.13265306122449
$ws.Range("J9").Value  = 0.142857142857143
$ws.Range("J10").Value = 0.13265306122449
$ws.Range("J11").Value = 0.13265306122449
